# MSDSCapstoneProjectDesignPlan.xlsx - "Updated as of 5/17" edit
#
# Summary of the change:
#  - 4 new/edited meeting rows are added right after "Meeting #15" (old row 41):
#       "Meeting # 16", "Meeting with Advisor (Pablo) # 3" (replaces the text of the
#       old "Meeting with Advisor #7" entry), "Meeting with Advisor # 7" (new),
#       "Meeting # 17" (new) - each completed on specific dates in May 2018.
#  - Everything that used to live at rows 43-56 shifts down by 4 rows (to 47-60) to
#    make room, with two blank spacer rows (45-46) in between (matching the old
#    spacer row 42).
#  - "Second Paper Draft (B)" (now at row 47) gets a real Start/End date and is
#    marked Completed instead of Not Started.
#  - The Gantt chart's series ranges grow to match the new data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: shift the block of rows 43-56 down to rows 47-60 (bottom-up so we
# never overwrite a row before it has been read).
# ---------------------------------------------------------------------------
for ($old = 56; $old -ge 43; $old--) {
    $new = $old + 4
    $ws.Range("B$old`:F$old").Copy()
    $ws.Range("B$new").PasteSpecial(-4104)   # xlPasteAll
}
$ws.Application.CutCopyMode = $false

# Clear the rows vacated by the shift that are *not* going to be overwritten
# below (43, 44, 45, 46 all get rewritten next, so nothing extra to clear here).

# ---------------------------------------------------------------------------
# Step 2: build the 4 new meeting rows (41-44), using row 40 ("Meeting #15",
# already Completed) as the formatting template so no new cell styles are
# created.
# ---------------------------------------------------------------------------
$ws.Range("B40:F40").Copy()
$ws.Range("B41:F44").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

function Set-MeetingRow($row, $name, $dateSerial) {
    $ws.Range("B$row").Value = $name
    $ws.Range("C$row").Value = $dateSerial
    $ws.Range("D$row").Value = $dateSerial
    $ws.Range("E$row").Formula = '=IF(ISBLANK(C' + $row + '),"", (D' + $row + '-C' + $row + '+1))'
    $ws.Range("F$row").Value = "Completed"
}

Set-MeetingRow 41 "Meeting # 16" 43223
Set-MeetingRow 42 "Meeting with Advisor (Pablo) # 3" 43229
Set-MeetingRow 43 "Meeting with Advisor # 7" 43235
Set-MeetingRow 44 "Meeting # 17" 43237

# ---------------------------------------------------------------------------
# Step 3: rows 45-46 are blank spacer rows (matching the old spacer row 42) -
# clear any content the shift may have left there and restore the plain
# formatting.
# ---------------------------------------------------------------------------
$ws.Range("B40:F40").Copy()
$ws.Range("B45:F46").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B45:F46").ClearContents()

# ---------------------------------------------------------------------------
# Step 4: "Second Paper Draft (B)" (now at row 47) is marked Completed with
# real start/end dates instead of the empty "Not Started" placeholder.
# ---------------------------------------------------------------------------
$ws.Range("C47").Value = 43132
$ws.Range("D47").Value = 43238
$ws.Range("F47").Value = "Completed"

# ---------------------------------------------------------------------------
# Step 5: grow the Gantt chart's series so the new rows show up on the chart.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$series1 = $chart.SeriesCollection(1)
$series1.Formula = "=SERIES(""Start Date"",'Basic Manual Gantt Chart'!`$B`$8:`$B`$53,'Basic Manual Gantt Chart'!`$C`$8:`$C`$51,1)"
$series2 = $chart.SeriesCollection(2)
$series2.Formula = "=SERIES(""Duration"",'Basic Manual Gantt Chart'!`$B`$8:`$B`$53,'Basic Manual Gantt Chart'!`$E`$8:`$E`$51,2)"

Write-Host "Edit applied"
